$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "studies" (sheet2): rename headers, drop/add a couple of
# data cells, add a new trailing column "reference_year".
# ---------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item("studies")
$wsStudies.Range("B1").Value = "study_label"
$wsStudies.Range("C1").Value = "description"
$wsStudies.Range("D1").Value = "access_level"
$wsStudies.Range("E1").Value = "contributors"
$wsStudies.Range("F1").Value = "reference"
$wsStudies.Range("G1").Value = "reference_year"
$wsStudies.Range("C2").ClearContents()
$wsStudies.Range("D2").Value = "public"
$wsStudies.Columns.Item(7).ColumnWidth = 20

# ---------------------------------------------------------------
# Sheet "surveys" (sheet3): rename first header, split the old
# "spatial_notes" column into "location_method"/"location_notes",
# and split a new "time_method" column in ahead of "time_notes".
# ---------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Range("A1").Value = "study_id"
$wsSurveys.Range("G1").Value = "location_method"

# Insert a new column at H: collection_start..time_notes shift right.
$wsSurveys.Columns.Item(8).Insert()
$wsSurveys.Range("H1").Value = "location_notes"
$wsSurveys.Range("H2").Value = "example data"
$wsSurveys.Range("G2").ClearContents()

# Insert a second new column before the (now shifted) time_notes column.
$wsSurveys.Columns.Item(12).Insert()
$wsSurveys.Range("L1").Value = "time_method"
$wsSurveys.Range("L2").Clear()

# ---------------------------------------------------------------
# Sheet "counts" (sheet4): rename first two headers only.
# ---------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("A1").Value = "study_id"
$wsCounts.Range("B1").Value = "survey_id"

# ---------------------------------------------------------------
# Update the selections / active cells on each sheet to match the
# edited workbook, finishing on "studies" so it becomes the
# workbook's active tab.
# ---------------------------------------------------------------
$wsCounts.Range("E8").Select() | Out-Null
$wsSurveys.Range("L2").Select() | Out-Null
$wsStudies.Range("D6").Select() | Out-Null
